# Auto-generated edit script: updates market-price / profit columns
# (currentAveragePrice*, LevePriceNQ/HQ, LeveProfitNQ/HQ) across several
# Leve-profit tracker sheets, matching a scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 9
$ws.Range("I16").Value = 9
$ws.Range("K16").Value = 9
$ws.Range("M16").Value = 221
$ws.Range("H19").Value = 1301.9
$ws.Range("I19").Value = 925.8889
$ws.Range("J19").Value = 1609.5454
$ws.Range("K19").Value = 925.8889
$ws.Range("L19").Value = 1609.5454
$ws.Range("M19").Value = -750.8889
$ws.Range("N19").Value = -1959.5454
$ws.Range("H28").Value = 863
$ws.Range("I28").Value = 784
$ws.Range("K28").Value = 784
$ws.Range("M28").Value = -299
$ws.Range("H40").Value = 4228.5
$ws.Range("I40").Value = 2818.1
$ws.Range("J40").Value = 5638.9
$ws.Range("K40").Value = 2818.1
$ws.Range("L40").Value = 5638.9
$ws.Range("M40").Value = -2643.1
$ws.Range("N40").Value = -5988.9
$ws.Range("H51").Value = 5998
$ws.Range("I51").Value = 5998
$ws.Range("K51").Value = 5998
$ws.Range("M51").Value = -5514
$ws.Range("H58").Value = 1330.5
$ws.Range("J58").Value = 1974.75
$ws.Range("L58").Value = 5924.25
$ws.Range("N58").Value = -6224.25
$ws.Range("H70").Value = 4774.125
$ws.Range("I70").Value = 4548.25
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 13644.75
$ws.Range("L70").Value = 15000
$ws.Range("M70").Value = -13374.75
$ws.Range("N70").Value = -15540
$ws.Range("H73").Value = 4774.125
$ws.Range("I73").Value = 4548.25
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 13644.75
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = -12708.75
$ws.Range("N73").Value = -16872
$ws.Range("H74").Value = 8527.308000000001
$ws.Range("I74").Value = 3979.5715
$ws.Range("K74").Value = 3979.5715
$ws.Range("M74").Value = -3043.5715
$ws.Range("H77").Value = 8527.308000000001
$ws.Range("I77").Value = 3979.5715
$ws.Range("K77").Value = 19897.8575
$ws.Range("M77").Value = -15217.8575
$ws.Range("H98").Value = 2561.8948
$ws.Range("J98").Value = 3356.3333
$ws.Range("L98").Value = 3356.3333
$ws.Range("N98").Value = -6352.3333
$ws.Range("H100").Value = 1787.6666
$ws.Range("I100").Value = 1348.1666
$ws.Range("J100").Value = 2666.6667
$ws.Range("K100").Value = 1348.1666
$ws.Range("L100").Value = 2666.6667
$ws.Range("M100").Value = -807.1666
$ws.Range("N100").Value = -3748.6667
$ws.Range("H103").Value = 854.1429000000001
$ws.Range("J103").Value = 887.5
$ws.Range("L103").Value = 2662.5
$ws.Range("N103").Value = -3834.5
$ws.Range("H122").Value = 2561.8948
$ws.Range("J122").Value = 3356.3333
$ws.Range("L122").Value = 10068.9999
$ws.Range("N122").Value = -14968.9999
$ws.Range("H131").Value = 68556.87
$ws.Range("I131").Value = 78996.766
$ws.Range("K131").Value = 236990.298
$ws.Range("M131").Value = -231950.298
$ws.Range("H132").Value = 993.7692
$ws.Range("I132").Value = 993.3913
$ws.Range("J132").Value = 996.6667
$ws.Range("K132").Value = 2980.1739
$ws.Range("L132").Value = 2990.0001
$ws.Range("M132").Value = -450.1738999999998
$ws.Range("N132").Value = -8050.0001
$ws.Range("H137").Value = 1736.56
$ws.Range("I137").Value = 1470.1052
$ws.Range("K137").Value = 4410.3156
$ws.Range("M137").Value = -1860.3156
$ws.Range("H138").Value = 3081.7144
$ws.Range("I138").Value = 2611.2942
$ws.Range("J138").Value = 3401.6
$ws.Range("K138").Value = 7833.882599999999
$ws.Range("L138").Value = 10204.8
$ws.Range("M138").Value = -2693.882599999999
$ws.Range("N138").Value = -20484.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4863.4062
$ws.Range("I32").Value = 4863.4062
$ws.Range("K32").Value = 4863.4062
$ws.Range("M32").Value = -4576.4062
$ws.Range("H61").Value = 7456.853
$ws.Range("I61").Value = 6211.48
$ws.Range("K61").Value = 6211.48
$ws.Range("M61").Value = -5999.48
$ws.Range("H74").Value = 2307.1462
$ws.Range("I74").Value = 1958.3334
$ws.Range("J74").Value = 3258.4546
$ws.Range("K74").Value = 1958.3334
$ws.Range("L74").Value = 3258.4546
$ws.Range("M74").Value = -1084.3334
$ws.Range("N74").Value = -5006.4546
$ws.Range("H77").Value = 2307.1462
$ws.Range("I77").Value = 1958.3334
$ws.Range("J77").Value = 3258.4546
$ws.Range("K77").Value = 9791.666999999999
$ws.Range("L77").Value = 16292.273
$ws.Range("M77").Value = -5423.666999999999
$ws.Range("N77").Value = -25028.273
$ws.Range("H102").Value = 4053.4736
$ws.Range("I102").Value = 3334.4666
$ws.Range("K102").Value = 3334.4666
$ws.Range("M102").Value = -1712.4666
$ws.Range("H122").Value = 4460
$ws.Range("I122").Value = 4460
$ws.Range("K122").Value = 13380
$ws.Range("M122").Value = -10930
$ws.Range("H132").Value = 4874.5884
$ws.Range("I132").Value = 3645.8147
$ws.Range("K132").Value = 10937.4441
$ws.Range("M132").Value = -8407.444100000001
$ws.Range("H136").Value = 7456.853
$ws.Range("I136").Value = 6211.48
$ws.Range("K136").Value = 18634.44
$ws.Range("M136").Value = -16084.44
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 40812.25
$ws.Range("I2").Value = 30000
$ws.Range("J2").Value = 44416.332
$ws.Range("K2").Value = 30000
$ws.Range("L2").Value = 44416.332
$ws.Range("M2").Value = -29887
$ws.Range("N2").Value = -44642.332
$ws.Range("H86").Value = 2231.6365
$ws.Range("I86").Value = 2181
$ws.Range("J86").Value = 2366.6667
$ws.Range("K86").Value = 2181
$ws.Range("L86").Value = 2366.6667
$ws.Range("M86").Value = -1058
$ws.Range("N86").Value = -4612.6667
$ws.Range("H89").Value = 2231.6365
$ws.Range("I89").Value = 2181
$ws.Range("J89").Value = 2366.6667
$ws.Range("K89").Value = 10905
$ws.Range("L89").Value = 11833.3335
$ws.Range("M89").Value = -5289
$ws.Range("N89").Value = -23065.3335
$ws.Range("H94").Value = 4874.5
$ws.Range("I94").Value = 4549.4
$ws.Range("J94").Value = 6500
$ws.Range("K94").Value = 4549.4
$ws.Range("L94").Value = 6500
$ws.Range("M94").Value = -4098.4
$ws.Range("N94").Value = -7402
$ws.Range("H99").Value = 2978.0605
$ws.Range("I99").Value = 2068.3157
$ws.Range("J99").Value = 4212.7144
$ws.Range("K99").Value = 2068.3157
$ws.Range("L99").Value = 4212.7144
$ws.Range("M99").Value = -570.3157000000001
$ws.Range("N99").Value = -7208.7144
$ws.Range("H107").Value = 2584.9614
$ws.Range("I107").Value = 2226.6956
$ws.Range("K107").Value = 2226.6956
$ws.Range("M107").Value = -306.6956
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 40723
$ws.Range("J28").Value = 40723
$ws.Range("L28").Value = 40723
$ws.Range("N28").Value = -41213
$ws.Range("H31").Value = 5290.9565
$ws.Range("I31").Value = 4105.933
$ws.Range("J31").Value = 7512.875
$ws.Range("K31").Value = 4105.933
$ws.Range("L31").Value = 7512.875
$ws.Range("M31").Value = -3810.933
$ws.Range("N31").Value = -8102.875
$ws.Range("H34").Value = 5290.9565
$ws.Range("I34").Value = 4105.933
$ws.Range("J34").Value = 7512.875
$ws.Range("K34").Value = 4105.933
$ws.Range("L34").Value = 7512.875
$ws.Range("M34").Value = -3903.933
$ws.Range("N34").Value = -7916.875
$ws.Range("H35").Value = 3201.3333
$ws.Range("I35").Value = 1830.2858
$ws.Range("J35").Value = 8000
$ws.Range("K35").Value = 1830.2858
$ws.Range("L35").Value = 8000
$ws.Range("M35").Value = -1536.2858
$ws.Range("N35").Value = -8588
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H58").Value = 6597.391
$ws.Range("I58").Value = 3031.6365
$ws.Range("J58").Value = 9866
$ws.Range("K58").Value = 3031.6365
$ws.Range("L58").Value = 9866
$ws.Range("M58").Value = -2828.6365
$ws.Range("N58").Value = -10272
$ws.Range("H86").Value = 9141.888999999999
$ws.Range("I86").Value = 8570
$ws.Range("J86").Value = 9599.4
$ws.Range("K86").Value = 8570
$ws.Range("L86").Value = 9599.4
$ws.Range("M86").Value = -7447
$ws.Range("N86").Value = -11845.4
$ws.Range("H89").Value = 9141.888999999999
$ws.Range("I89").Value = 8570
$ws.Range("J89").Value = 9599.4
$ws.Range("K89").Value = 42850
$ws.Range("L89").Value = 47997
$ws.Range("M89").Value = -37234
$ws.Range("N89").Value = -59229
$ws.Range("H97").Value = 34989.8
$ws.Range("J97").Value = 36237.25
$ws.Range("L97").Value = 36237.25
$ws.Range("N97").Value = -38219.25
$ws.Range("H105").Value = 2165
$ws.Range("I105").Value = 2536.6667
$ws.Range("K105").Value = 2536.6667
$ws.Range("M105").Value = -789.6667000000002
$ws.Range("H108").Value = 40000
$ws.Range("J108").Value = 40000
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H111").Value = 77919
$ws.Range("J111").Value = 77919
$ws.Range("L111").Value = 77919
$ws.Range("N111").Value = -86099
$ws.Range("H112").Value = 75973
$ws.Range("J112").Value = 75973
$ws.Range("L112").Value = 75973
$ws.Range("N112").Value = -78927
$ws.Range("H132").Value = 4401.933
$ws.Range("I132").Value = 2448.25
$ws.Range("J132").Value = 6634.7144
$ws.Range("K132").Value = 7344.75
$ws.Range("L132").Value = 19904.1432
$ws.Range("M132").Value = -4814.75
$ws.Range("N132").Value = -24964.1432
$ws.Range("H134").Value = 5824.5
$ws.Range("I134").Value = 4088.2778
$ws.Range("J134").Value = 11033.167
$ws.Range("K134").Value = 12264.8334
$ws.Range("L134").Value = 33099.501
$ws.Range("M134").Value = -9729.8334
$ws.Range("N134").Value = -38169.501
$ws.Range("H136").Value = 6597.391
$ws.Range("I136").Value = 3031.6365
$ws.Range("J136").Value = 9866
$ws.Range("K136").Value = 9094.9095
$ws.Range("L136").Value = 29598
$ws.Range("M136").Value = -6544.9095
$ws.Range("N136").Value = -34698
$ws.Range("H141").Value = 36633.332
$ws.Range("J141").Value = 39900
$ws.Range("L141").Value = 39900
$ws.Range("N141").Value = -50260
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1250
$ws.Range("H76").Value = 4000
$ws.Range("I76").Value = 4000
$ws.Range("K76").Value = 12000
$ws.Range("M76").Value = -11617
$ws.Range("H79").Value = 4000
$ws.Range("I79").Value = 4000
$ws.Range("K79").Value = 12000
$ws.Range("M79").Value = -10674
$ws.Range("H104").Value = 4675.3335
$ws.Range("I104").Value = 8026
$ws.Range("J104").Value = 3000
$ws.Range("K104").Value = 24078
$ws.Range("L104").Value = 9000
$ws.Range("M104").Value = -21457
$ws.Range("N104").Value = -14242
$ws.Range("H121").Value = 16667048
$ws.Range("J121").Value = 25000356
$ws.Range("L121").Value = 75001068
$ws.Range("N121").Value = -75003688
$ws.Range("H128").Value = 521830.78
$ws.Range("I128").Value = 521830.78
$ws.Range("K128").Value = 1565492.34
$ws.Range("M128").Value = -1560512.34
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 500
$ws.Range("J10").Value = 500
$ws.Range("L10").Value = 500
$ws.Range("N10").Value = -838
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H99").Value = 11715.714
$ws.Range("I99").Value = 7335
$ws.Range("K99").Value = 7335
$ws.Range("M99").Value = -5089
$ws.Range("H107").Value = 579.375
$ws.Range("I107").Value = 684.625
$ws.Range("J107").Value = 368.875
$ws.Range("K107").Value = 684.625
$ws.Range("L107").Value = 368.875
$ws.Range("M107").Value = 1235.375
$ws.Range("N107").Value = -4208.875
$ws.Range("H113").Value = 445648.78
$ws.Range("J113").Value = 549.5
$ws.Range("L113").Value = 549.5
$ws.Range("N113").Value = -4889.5
$ws.Range("H122").Value = 1597.1428
$ws.Range("I122").Value = 1480.8334
$ws.Range("K122").Value = 4442.5002
$ws.Range("M122").Value = -1992.5002
$ws.Range("H123").Value = 39072.418
$ws.Range("J123").Value = 39072.418
$ws.Range("L123").Value = 39072.418
$ws.Range("N123").Value = -43972.418
$ws.Range("H132").Value = 2956
$ws.Range("I132").Value = 2744.6956
$ws.Range("J132").Value = 3361
$ws.Range("K132").Value = 8234.086800000001
$ws.Range("L132").Value = 10083
$ws.Range("M132").Value = -5704.086800000001
$ws.Range("N132").Value = -15143
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H22").Value = 3333.1667
$ws.Range("J22").Value = 3699.8
$ws.Range("L22").Value = 3699.8
$ws.Range("N22").Value = -4289.8
$ws.Range("H25").Value = 28832.666
$ws.Range("J25").Value = 29999
$ws.Range("L25").Value = 29999
$ws.Range("N25").Value = -30459
$ws.Range("H27").Value = 3333.1667
$ws.Range("J27").Value = 3699.8
$ws.Range("L27").Value = 3699.8
$ws.Range("N27").Value = -3913.8
$ws.Range("H41").Value = 1600
$ws.Range("I41").Value = 1600
$ws.Range("K41").Value = 1600
$ws.Range("M41").Value = -1162
$ws.Range("H46").Value = 11161.692
$ws.Range("J46").Value = 12554.728
$ws.Range("L46").Value = 12554.728
$ws.Range("N46").Value = -12930.728
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H55").Value = 1126.8334
$ws.Range("J55").Value = 758.3333
$ws.Range("L55").Value = 758.3333
$ws.Range("N55").Value = -1104.3333
$ws.Range("H61").Value = 48774.59
$ws.Range("I61").Value = 62090.65
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 62090.65
$ws.Range("L61").Value = 3500
$ws.Range("M61").Value = -61888.65
$ws.Range("N61").Value = -3904
$ws.Range("H68").Value = 12189.036
$ws.Range("J68").Value = 10873.25
$ws.Range("L68").Value = 10873.25
$ws.Range("N68").Value = -12371.25
$ws.Range("H71").Value = 12189.036
$ws.Range("J71").Value = 10873.25
$ws.Range("L71").Value = 54366.25
$ws.Range("N71").Value = -61854.25
$ws.Range("H93").Value = 19205.883
$ws.Range("I93").Value = 1100
$ws.Range("J93").Value = 31880
$ws.Range("K93").Value = 1100
$ws.Range("L93").Value = 31880
$ws.Range("M93").Value = 148
$ws.Range("N93").Value = -34376
$ws.Range("H113").Value = 48774.59
$ws.Range("I113").Value = 62090.65
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 62090.65
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = -59920.65
$ws.Range("N113").Value = -7840
$ws.Range("H136").Value = 7876.2
$ws.Range("J136").Value = 377
$ws.Range("L136").Value = 1131
$ws.Range("N136").Value = -6231
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 40303.75
$ws.Range("J70").Value = 40303.75
$ws.Range("L70").Value = 40303.75
$ws.Range("N70").Value = -40933.75
$ws.Range("H73").Value = 40303.75
$ws.Range("J73").Value = 40303.75
$ws.Range("L73").Value = 40303.75
$ws.Range("N73").Value = -42487.75
$ws.Range("H81").Value = 1599.25
$ws.Range("I81").Value = 1599.25
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3198.5
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -2137.5
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 1599.25
$ws.Range("I84").Value = 1599.25
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 15992.5
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -10688.5
$ws.Range("N84").ClearContents()
$ws.Range("H92").Value = 77000
$ws.Range("J92").Value = 77000
$ws.Range("L92").Value = 77000
$ws.Range("N92").Value = -81992
$ws.Range("H107").Value = 847.6667
$ws.Range("I107").Value = 733.2778
$ws.Range("K107").Value = 2199.8334
$ws.Range("M107").Value = -279.8334
$ws.Range("H122").Value = 3576.9524
$ws.Range("I122").Value = 2031.4445
$ws.Range("J122").Value = 4736.0835
$ws.Range("K122").Value = 6094.333500000001
$ws.Range("L122").Value = 14208.2505
$ws.Range("M122").Value = -3644.333500000001
$ws.Range("N122").Value = -19108.2505
$ws.Range("H132").Value = 4159.731
$ws.Range("I132").Value = 3299.6667
$ws.Range("K132").Value = 9899.000100000001
$ws.Range("M132").Value = -7369.000100000001
$ws.Range("H133").Value = 111163.75
$ws.Range("J133").Value = 111163.75
$ws.Range("L133").Value = 111163.75
$ws.Range("N133").Value = -121283.75
$ws.Range("H136").Value = 6470.2607
$ws.Range("I136").Value = 4368.25
$ws.Range("J136").Value = 8763.362999999999
$ws.Range("K136").Value = 13104.75
$ws.Range("L136").Value = 26290.089
$ws.Range("M136").Value = -10554.75
$ws.Range("N136").Value = -31390.089
